{"js": "// Extend the final bullet paragraph (\"Implementaremos un alfabeto como gu\u00eda\n// para balancear (ordenar) el grafo.\") with the new clause about comparing\n// names alphabetically and avoiding the usual \"<\"/\">\" operators, and move\n// the (hidden) \"_GoBack\" bookmark from the earlier \"Eliminar ... carpeta.\"\n// paragraph to right before the final period of the last paragraph.\n\nconst body = context.document.body;\n\n// The \"_GoBack\" bookmark currently sits inside the \"Eliminar ruta, archivo y\n// carpeta.\" paragraph. Drop it there; it will be re-created at the new\n// location below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Locate the unique, unmodified beginning of the sentence we need to extend.\nconst oldCore = \"Implementaremos un alfabeto como gu\u00eda para balancear (ordenar) el grafo\";\nconst matches = body.search(oldCore, { matchCase: true, matchWholeWord: false });\nmatches.load(\"text\");\nawait context.sync();\n\nif (matches.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one match for \"${oldCore}\", found ${matches.items.length}`\n  );\n}\n\nconst target = matches.items[0];\n\n// New text to splice in right after \"...el grafo\" and before the trailing\n// period that already ends the paragraph.\nconst addition =\n  \", se comparara cada uno de los nombres de los enlaces, archivos y carpetas para ordenarlos alfab\u00e9ticamente; \" +\n  \"evitando usar los operadores l\u00f3gicos usuales \\u201c<\\u201d y \\u201c>\\u201d para las cadenas\";\n\nconst insertedRange = target.insertText(addition, Word.InsertLocation.after);\nawait context.sync();\n\n// Re-create \"_GoBack\" as a collapsed bookmark right after the newly inserted\n// text, i.e. immediately before the paragraph's final \".\".\nconst endOfInsertion = insertedRange.getRange(Word.RangeLocation.end);\nendOfInsertion.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Extend the final bullet paragraph (\"Implementaremos un alfabeto como gu\u00eda\n# para balancear (ordenar) el grafo.\") with the new clause about comparing\n# names alphabetically and avoiding the usual \"<\"/\">\" operators, and move\n# the (hidden) \"_GoBack\" bookmark from the earlier \"Eliminar ... carpeta.\"\n# paragraph to right before the final period of the last paragraph.\n\n$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark currently sits inside the \"Eliminar ruta, archivo y\n# carpeta.\" paragraph. Drop it there; it will be re-created at the new\n# location below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the unique, unmodified beginning of the sentence we need to extend\n# (everything up to, but excluding, the trailing period).\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"Implementaremos un alfabeto como gu\u00eda para balancear (ordenar) el grafo\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target sentence to extend.\"\n}\n\n# New text to splice in right after \"...el grafo\" and before the trailing\n# period that already ends the paragraph.\n$addition = \", se comparara cada uno de los nombres de los enlaces, archivos y carpetas para ordenarlos alfab\u00e9ticamente; evitando usar los operadores l\u00f3gicos usuales \u201c<\u201d y \u201c>\u201d para las cadenas\"\n\n$insertPoint = $searchRange.Duplicate\n$insertPoint.Collapse(0)            # wdCollapseEnd\n$insertPoint.InsertAfter($addition)\n\n# Re-create \"_GoBack\" as a collapsed bookmark right after the newly inserted\n# text, i.e. immediately before the paragraph's final \".\".\n$insertPoint.Collapse(0)            # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n"}
